# "Changed Radius Values, Skidpad Wierd Output"
#
# Updates the radius/section values on the Slalom and UTurn sheets, updates
# the Skidpad sheet's values, and leaves the Skidpad tab as the active /
# selected sheet (instead of Legend) so the user can look at its output.

$wb = $excel.ActiveWorkbook

# Keep Legend's own selection as-is (C6) while we're poking at the other
# sheets - it just won't be the active tab anymore once we activate Skidpad
# last.
$wsLegend = $wb.Worksheets.Item("Legend")
$wsLegend.Activate()
$wsLegend.Range("C6").Select()

# UTurn: single row of values changes, and a cell selection (B2) is recorded
# for the first time.
$wsUTurn = $wb.Worksheets.Item("UTurn")
$wsUTurn.Activate()
$wsUTurn.Range("B1").Value = 13
$wsUTurn.Range("C1").Value = 4
$wsUTurn.Range("B2").Select()

# Slalom: all five rows get new radius/length values, selection moves to B2.
$wsSlalom = $wb.Worksheets.Item("Slalom")
$wsSlalom.Activate()
for ($r = 1; $r -le 5; $r++) {
    $wsSlalom.Cells.Item($r, 2).Value = 10
    $wsSlalom.Cells.Item($r, 3).Value = 8.5
}
$wsSlalom.Range("B2").Select()

# Skidpad: both rows get new values, selection moves to B3, and this sheet
# becomes the active/selected tab (activeTab=2 / tabSelected=1).
$wsSkidpad = $wb.Worksheets.Item("Skidpad")
$wsSkidpad.Activate()
$wsSkidpad.Range("B1").Value = 58
$wsSkidpad.Range("C1").Value = 9.125
$wsSkidpad.Range("B2").Value = 58
$wsSkidpad.Range("C2").Value = 9.125
$wsSkidpad.Range("B3").Select()
